$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Unique" column (L) mirroring the existing "Multivalued" column (K):
# header text in row 4, and literal text "FALSE" values in rows 5-7.
# A leading apostrophe forces the "FALSE" entries to be stored as text
# (matching column K) instead of being auto-coerced to a boolean.
$ws.Range("L4").Value = "Unique"
$ws.Range("L5").Value = "'FALSE"
$ws.Range("L6").Value = "'FALSE"
$ws.Range("L7").Value = "'FALSE"

# Copy column K's formatting (bold header style / boolean-style text format)
# onto the new column L so it looks just like K.
$ws.Range("K4:K7").Copy()
$ws.Range("L4:L7").PasteSpecial(-4122) # xlPasteFormats

# Move the active selection to the new column, matching the author's edit.
$ws.Range("L4:L7").Select()
